$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H33").Value = 342.9375
$ws.Range("I33").Value = 177.46153
$ws.Range("J33").Value = 1060
$ws.Range("K33").Value = 177.46153
$ws.Range("L33").Value = 1060
$ws.Range("M33").Value = 51.53846999999999
$ws.Range("N33").Value = -1518
$ws.Range("H53").Value = 172.41176
$ws.Range("I53").Value = 112.416664
$ws.Range("J53").Value = 316.4
$ws.Range("K53").Value = 112.416664
$ws.Range("L53").Value = 316.4
$ws.Range("M53").Value = 524.583336
$ws.Range("N53").Value = -1590.4
$ws.Range("H81").Value = 35000
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -36996
$ws.Range("H84").Value = 35000
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -114984
$ws.Range("H118").Value = 38769750
$ws.Range("I118").Value = 42000316
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 126000948
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = -125999291
$ws.Range("N118").Value = -12314
$ws.Range("H129").Value = 1229.0938
$ws.Range("I129").Value = 535.4545000000001
$ws.Range("J129").Value = 1592.4286
$ws.Range("K129").Value = 1606.3635
$ws.Range("L129").Value = 4777.2858
$ws.Range("M129").Value = 3393.6365
$ws.Range("N129").Value = -14777.2858
$ws.Range("H137").Value = 4546986.5
$ws.Range("I137").Value = 1786841.4
$ws.Range("J137").Value = 20003800
$ws.Range("K137").Value = 5360524.199999999
$ws.Range("L137").Value = 60011400
$ws.Range("M137").Value = -5357974.199999999
$ws.Range("N137").Value = -60016500
$ws.Range("H138").Value = 2651.7046
$ws.Range("J138").Value = 2414.7437
$ws.Range("L138").Value = 7244.2311
$ws.Range("N138").Value = -17524.2311
$ws.Range("H139").Value = 34850
$ws.Range("J139").Value = 39800
$ws.Range("L139").Value = 39800
$ws.Range("N139").Value = -50080
$ws.Range("H141").Value = 3450
$ws.Range("I141").Value = 3195
$ws.Range("K141").Value = 9585
$ws.Range("M141").Value = -4405

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H5").Value = 970
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H122").Value = 1663.9445
$ws.Range("I122").Value = 1395.8
$ws.Range("J122").Value = 3004.6667
$ws.Range("K122").Value = 4187.4
$ws.Range("L122").Value = 9014.000100000001
$ws.Range("M122").Value = -1737.4
$ws.Range("N122").Value = -13914.0001

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H4").Value = 970
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H59").Value = 64950
$ws.Range("J59").Value = 64950
$ws.Range("L59").Value = 64950
$ws.Range("N59").Value = -66644
$ws.Range("H138").Value = 39800
$ws.Range("J138").Value = 39800
$ws.Range("L138").Value = 39800
$ws.Range("N138").Value = -50080

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H132").Value = 4204.5835
$ws.Range("I132").Value = 3814
$ws.Range("J132").Value = 4985.75
$ws.Range("K132").Value = 11442
$ws.Range("L132").Value = 14957.25
$ws.Range("M132").Value = -8912
$ws.Range("N132").Value = -20017.25
$ws.Range("H134").Value = 6669.609
$ws.Range("I134").Value = 10157.615
$ws.Range("J134").Value = 2135.2
$ws.Range("K134").Value = 30472.845
$ws.Range("L134").Value = 6405.599999999999
$ws.Range("M134").Value = -27937.845
$ws.Range("N134").Value = -11475.6

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H121").Value = 635704.6
$ws.Range("I121").Value = 328.57144
$ws.Range("J121").Value = 953392.6
$ws.Range("K121").Value = 985.71432
$ws.Range("L121").Value = 2860177.8
$ws.Range("M121").Value = 324.28568
$ws.Range("N121").Value = -2862797.8
$ws.Range("H131").Value = 905.25
$ws.Range("I131").Value = 345.75
$ws.Range("J131").Value = 940.21875
$ws.Range("K131").Value = 1037.25
$ws.Range("L131").Value = 2820.65625
$ws.Range("M131").Value = 4002.75
$ws.Range("N131").Value = -12900.65625

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H82").Value = 44000
$ws.Range("J82").Value = 44000
$ws.Range("L82").Value = 44000
$ws.Range("N82").Value = -44766
$ws.Range("H85").Value = 44000
$ws.Range("J85").Value = 44000
$ws.Range("L85").Value = 44000
$ws.Range("N85").Value = -46652
$ws.Range("H101").Value = 42888.668
$ws.Range("J101").Value = 42888.668
$ws.Range("L101").Value = 42888.668
$ws.Range("N101").Value = -49378.668
$ws.Range("H116").Value = 31000
$ws.Range("J116").Value = 31000
$ws.Range("L116").Value = 31000
$ws.Range("N116").Value = -40178
$ws.Range("H122").Value = 1688.1562
$ws.Range("I122").Value = 1805.9048
$ws.Range("J122").Value = 1463.3636
$ws.Range("K122").Value = 5417.7144
$ws.Range("L122").Value = 4390.0908
$ws.Range("M122").Value = -2967.7144
$ws.Range("N122").Value = -9290.0908
$ws.Range("H132").Value = 1670.4517
$ws.Range("I132").Value = 1349.5
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 4048.5
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -1518.5
$ws.Range("N132").Value = -19058

# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 3706095
$ws.Range("I7").Value = 5884001
$ws.Range("J7").Value = 3654.6
$ws.Range("K7").Value = 5884001
$ws.Range("L7").Value = 3654.6
$ws.Range("M7").Value = -5883889
$ws.Range("N7").Value = -3878.6
$ws.Range("H22").Value = 913.97675
$ws.Range("I22").Value = 659.8333
$ws.Range("J22").Value = 1012.35486
$ws.Range("K22").Value = 659.8333
$ws.Range("L22").Value = 1012.35486
$ws.Range("M22").Value = -364.8333
$ws.Range("N22").Value = -1602.35486
$ws.Range("H27").Value = 913.97675
$ws.Range("I27").Value = 659.8333
$ws.Range("J27").Value = 1012.35486
$ws.Range("K27").Value = 659.8333
$ws.Range("L27").Value = 1012.35486
$ws.Range("M27").Value = -552.8333
$ws.Range("N27").Value = -1226.35486
$ws.Range("H46").Value = 1198.9474
$ws.Range("I46").Value = 988.6667
$ws.Range("J46").Value = 1388.2
$ws.Range("K46").Value = 988.6667
$ws.Range("L46").Value = 1388.2
$ws.Range("M46").Value = -800.6667
$ws.Range("N46").Value = -1764.2
$ws.Range("H55").Value = 1024.8334
$ws.Range("I55").Value = 549.8333
$ws.Range("J55").Value = 1262.3334
$ws.Range("K55").Value = 549.8333
$ws.Range("L55").Value = 1262.3334
$ws.Range("M55").Value = -376.8333
$ws.Range("N55").Value = -1608.3334
$ws.Range("H126").Value = 3706095
$ws.Range("I126").Value = 5884001
$ws.Range("J126").Value = 3654.6
$ws.Range("K126").Value = 17652003
$ws.Range("L126").Value = 10963.8
$ws.Range("M126").Value = -17649533
$ws.Range("N126").Value = -15903.8
$ws.Range("H136").Value = 1558.8846
$ws.Range("I136").Value = 1272.1578
$ws.Range("J136").Value = 2337.1428
$ws.Range("K136").Value = 3816.4734
$ws.Range("L136").Value = 7011.428400000001
$ws.Range("M136").Value = -1266.4734
$ws.Range("N136").Value = -12111.4284

# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 3692.3076
$ws.Range("I132").Value = 3339
$ws.Range("J132").Value = 4487.25
$ws.Range("K132").Value = 10017
$ws.Range("L132").Value = 13461.75
$ws.Range("M132").Value = -7487
$ws.Range("N132").Value = -18521.75
$ws.Range("H136").Value = 1404.5952
$ws.Range("I136").Value = 1219.8055
$ws.Range("J136").Value = 2513.3333
$ws.Range("K136").Value = 3659.4165
$ws.Range("L136").Value = 7539.999899999999
$ws.Range("M136").Value = -1109.4165
$ws.Range("N136").Value = -12639.9999
